$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-Alert Template Import")

$ws.Range("A3").Value = "'JSSO1000244"
$ws.Range("B3").Value = "'JSSO1000244"
$ws.Range("C3").Value = "'JSSO1000244"
$ws.Range("AJ3").Value = "JSCN1000244"
$ws.Range("AN3").Value = "JSCN1000244"
